# Updated line text for the ChangingLines component.
# - Fixes the misspelling "Febuary" -> "February" for all DATE entries
#   that still had the typo (rows 275-316 in the "Log" table).
# - Fills in the remaining rows of the Daily Activits log (rows 317-320)
#   with the new entries for Feb 22-24 2018.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daily Activits")

# ---------------------------------------------------------------------
# 1. Correct the "Febuary" -> "February" typo in the DATE column (B)
#    for every date cell that still has the old spelling.
# ---------------------------------------------------------------------
$dateFixRows = 275,276,278,280,283,285,288,292,295,297,301,302,303,305,308,311,313

foreach ($r in $dateFixRows) {
    $cell = $ws.Cells.Item($r, 2)
    $old = [string]$cell.Value2
    $new = $old -replace "^Febuary", "February"
    $cell.Value = $new
}

# ---------------------------------------------------------------------
# 2. Fill in the last four rows of the Log table (317-320) with the new
#    activity entries.
# ---------------------------------------------------------------------

# Row 317: February 22 2018
$ws.Range("B317").Value = "February 22 2018"
$ws.Range("C317").Value = 0.75
$ws.Range("D317").Value = 0.77083333333333337
$ws.Range("F317").Value = "Changed add journal button's color. Added showing hexagram detail modal feature to search reading page."

# Row 318: February 23 2018
$ws.Range("B318").Value = "February 23 2018"
$ws.Range("C318").Value = 0.72916666666666663
$ws.Range("D318").Value = 0.75
$ws.Range("F318").Value = "Starting to extract some code relates to show hexagram detail modal."

# Row 319: (no date, continuation of Feb 23 work)
$ws.Range("C319").Value = 0.86111111111111116
$ws.Range("D319").Value = 1.0326388888888889
$ws.Range("F319").Value = "Finished refactoring for HexagramDetailModal."

# Row 320: February 24 2018 (only a start time so far)
$ws.Range("B320").Value = "February 24 2018"
$ws.Range("C320").Value = 0.65763888888888888

# ---------------------------------------------------------------------
# 3. Match the reviewer's on-screen position/selection from the commit
#    (scrolled down to the new rows, cell C320 selected).
# ---------------------------------------------------------------------
$ws.Range("C320").Select()
